$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update to latest Nomis data -------------------------------------------
# The "LatestPeriod" column (B) for the employment / self-employment /
# unemployment / inactivity rate + volume rows was "Jan-Dec 2022 data";
# refresh it to the newer Nomis release period.
$ws.Range("B2:B9").Value = "Apr 2022 - Mar 2023 data"

# --- Leave behind the scroll position / selection from the edit session ----
$win = $excel.ActiveWindow
$ws.Range("B2:B9").Select()
$win.ScrollRow = 17
$win.ScrollColumn = 1
